# COREESG_holdings.xlsx — refresh the "as of" date in the confidentiality
# footer and update the model's Weight / Percent Change figures for the
# NULG/NULV/NUMG/NUMV/NUSC/Total rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to edit, re-protect when done.
$ws.Unprotect()

# --- Footer text: bump the "as of" date from 2021-05-04 to 2021-05-05 ---
$footer = $ws.Range("A10").Value2
$ws.Range("A10").Value = $footer -replace "2021-05-04", "2021-05-05"

# --- Weight (D) / Percent Change (E) refresh for rows 2-7 ---
$ws.Range("D2").Value = 0.246969748564624
$ws.Range("E2").Value = -0.003161397670549237

$ws.Range("D3").Value = 0.496862446843175
$ws.Range("E3").Value = 0.002909283258397233

$ws.Range("D4").Value = 0.09761474137912279
$ws.Range("E4").Value = -0.001363459291001345

$ws.Range("D5").Value = 0.1016097926389595
$ws.Range("E5").Value = 0.004360861270100935

$ws.Range("D6").Value = 0.0569432705741186
$ws.Range("E6").Value = -0.002250731487733582

$ws.Range("D7").Value = 0.9999999999999999
$ws.Range("E7").Value = 0.0008465924817331327

# Restore sheet protection.
$ws.Protect()
